# "I0 and IF added"
# Adds two new columns (I = "I0", J = "IF") to the sheet: a header in row 1
# plus 66 data rows (rows 2-67), and extends the used range to A1:J67.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the formatting already used by the rest of the header row (H1):
# bold font, thin border, centered alignment.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Re-assert the values (PasteSpecial only touches formats, but make sure).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-67 -----------------------------------------------------
$iValues = @(9,7,7,6,8,5,7,6,6,6,9,5,5,8,7,7,9,7,6,5,7,9,5,6,6,6,10,9,9,9,7,6,6,7,8,6,7,6,7,8,6,9,4,7,7,8,7,6,8,7,8,5,7,6,5,9,7,9,9,7,6,7,6,7,5,4)
$jValues = @(9,7,8,6,8,5,7,6,6,6,9,5,6,8,7,7,9,7,6,6,7,9,6,6,6,6,10,9,9,9,7,7,6,7,8,6,7,6,8,8,6,9,4,7,7,8,7,7,8,7,8,5,8,7,6,9,8,9,9,7,6,7,6,7,5,4)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}

Write-Host "I0/IF columns added; dimension now A1:J67"
